# Preset Location Values Update
# - BasicInfo (sheet1): drop the three middle sample rows (old rows 3-5),
#   refresh the remaining two data rows with the new preset values, widen
#   column R, and move the selection/scroll position.
# - UpdateOptions (sheet5): move the selection only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BasicInfo
# ---------------------------------------------------------------------------
$basic = $wb.Worksheets.Item("BasicInfo")

# Remove the three obsolete sample rows; row 6 slides up to become row 3.
$basic.Range("A3:A5").EntireRow.Delete()

# Row 2 — refresh the preset account / reference / language / vendor values.
$basic.Range("A2").Value = "SA Test DRS"
$basic.Range("B2").Value = "21062021_1"
$basic.Range("Q2").Value = "English-Canada"
$basic.Range("R2").Value = "Apple,Bing,Facebook,Foursquare,Google,HERE,Tom Tom,Trip Advisor,Zomato"
$basic.Range("S2").Value = "20,2,4,3,1,10,30,18,29"

# Row 3 (previously row 6) — same preset-account refresh; also clear the
# one-off style that used to mark it as the last/special row.
$basic.Range("A3").Value = "SA Test DRS"
$basic.Range("A3").Style = "Normal"
$basic.Range("B3").Value = "21062021_2"
$basic.Range("Q3").Value = "English-Canada"

# Widen the vendor-list column to fit the longer preset text.
$basic.Columns("R").ColumnWidth = 77.67

# Move the active selection / scroll position.
$basic.Range("R9").Select()
$excel.ActiveWindow.ScrollColumn = 17
$excel.ActiveWindow.ScrollRow = 1

# ---------------------------------------------------------------------------
# UpdateOptions
# ---------------------------------------------------------------------------
$updateOptions = $wb.Worksheets.Item("UpdateOptions")
$updateOptions.Range("E10").Select()

# Re-activate BasicInfo so it stays the tab that's shown/selected.
$basic.Activate()
